# Daily refresh of the cryptos list (Coin / Link / Price / Volume(1h)).
# Updates Price (D) and Volume(1h) (E) figures for most rows, and for two
# row-pairs (24/25 and 46/47) the two coins swapped rank so their entire
# row contents (Coin, Link, Price, Volume) are exchanged.
#
# Some new Price values look like plain decimal numbers (e.g. "602.86",
# "33.00", "1.00"). Excel's COM layer auto-converts such literals to
# numbers, which would silently drop meaningful trailing zeros / add
# float noise. To keep column D as plain text (matching the sheet's
# original inline-string cells with no explicit number format), those
# values are written with a leading apostrophe (forces text entry) and
# then the cell style is reset back to "Normal" so no stray number
# format/style is left applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.516.66"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.497.51"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'602.86"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "'194.10"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -4.59%  "
$ws.Range("D10").Value = "'0.645"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "'53.03"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "'0.0000298"
$ws.Range("D13").Value = "'9.44"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "4.057.30"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'593.64"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "69.697.48"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "'12.63"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "3.497.39"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "'18.19"
$ws.Range("E22").Value = "  +6.91%  "
$ws.Range("D23").Value = "'5.24"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'101.43"
$ws.Range("E24").Value = "  -4.37%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'4.62"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "'3.13"
$ws.Range("E26").Value = "  +4.74%  "
$ws.Range("D27").Value = "'10.77"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "'33.00"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'4.30"
$ws.Range("E30").Value = "  +8.93%  "
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").Value = "'12.32"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'63.07"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "3.729.02"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "'3.10"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "0.0₃0810"
$ws.Range("E37").Value = "  +5.83%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "'0.388"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "'36.16"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "'491.02"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").Value = "'0.133"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.80"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.29"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("E51").Value = "  +10.06%  "

# Reset number format/style on cells that needed a text-prefix, to avoid leftover styling
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
